$wb = $excel.ActiveWorkbook

# Sheet ALC, row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1317.1389
$ws.Range("J17").Value = 1371.0883
$ws.Range("L17").Value = 4113.2649
$ws.Range("N17").Value = -4449.2649

# Sheet ALC, row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 19768498
$ws.Range("I116").Value = 19768498
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 19768498
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -19765056
$ws.Range("N116").ClearContents()

# Sheet ALC, row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 432085.9
$ws.Range("I132").Value = 506478.53
$ws.Range("J132").Value = 75001.2
$ws.Range("K132").Value = 1519435.59
$ws.Range("L132").Value = 225003.6
$ws.Range("M132").Value = -1516905.59
$ws.Range("N132").Value = -230063.6

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 11910602
$ws.Range("I138").Value = 7765.6665
$ws.Range("J138").Value = 13894408
$ws.Range("K138").Value = 23296.9995
$ws.Range("L138").Value = 41683224
$ws.Range("M138").Value = -18156.9995
$ws.Range("N138").Value = -41693504

# Sheet ARM, row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6598.476
$ws.Range("I2").Value = 8801.866
$ws.Range("J2").Value = 1090
$ws.Range("K2").Value = 8801.866
$ws.Range("L2").Value = 1090
$ws.Range("M2").Value = -8688.866
$ws.Range("N2").Value = -1316

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2287.0193
$ws.Range("I32").Value = 1764.2766
$ws.Range("K32").Value = 1764.2766
$ws.Range("M32").Value = -1477.2766

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1108
$ws.Range("I74").Value = 1250.909
$ws.Range("K74").Value = 1250.909
$ws.Range("M74").Value = -376.9090000000001

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1108
$ws.Range("I77").Value = 1250.909
$ws.Range("K77").Value = 6254.545
$ws.Range("M77").Value = -1886.545

# Sheet ARM, row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 6598.476
$ws.Range("I116").Value = 8801.866
$ws.Range("J116").Value = 1090
$ws.Range("K116").Value = 8801.866
$ws.Range("L116").Value = 1090
$ws.Range("M116").Value = -6507.866
$ws.Range("N116").Value = -5678

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3070.0857
$ws.Range("I132").Value = 2762.0476
$ws.Range("J132").Value = 3532.1428
$ws.Range("K132").Value = 8286.1428
$ws.Range("L132").Value = 10596.4284
$ws.Range("M132").Value = -5756.1428
$ws.Range("N132").Value = -15656.4284

# Sheet BSM, row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6598.476
$ws.Range("I3").Value = 8801.866
$ws.Range("J3").Value = 1090
$ws.Range("K3").Value = 8801.866
$ws.Range("L3").Value = 1090
$ws.Range("M3").Value = -8687.866
$ws.Range("N3").Value = -1318

# Sheet BSM, row 7
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 5333.3335
$ws.Range("I7").Value = 5000
$ws.Range("K7").Value = 5000
$ws.Range("M7").Value = -4887

# Sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2710.2144
$ws.Range("I134").Value = 1401.2
$ws.Range("J134").Value = 5982.75
$ws.Range("K134").Value = 4203.6
$ws.Range("L134").Value = 17948.25
$ws.Range("M134").Value = -1668.6
$ws.Range("N134").Value = -23018.25

# Sheet CRP, row 11
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 50000
$ws.Range("I11").Value = 50000
$ws.Range("K11").Value = 50000
$ws.Range("M11").Value = -49860

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1844.1951
$ws.Range("I31").Value = 1236.4615
$ws.Range("J31").Value = 2126.3572
$ws.Range("K31").Value = 1236.4615
$ws.Range("L31").Value = 2126.3572
$ws.Range("M31").Value = -941.4614999999999
$ws.Range("N31").Value = -2716.3572

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1844.1951
$ws.Range("I34").Value = 1236.4615
$ws.Range("J34").Value = 2126.3572
$ws.Range("K34").Value = 1236.4615
$ws.Range("L34").Value = 2126.3572
$ws.Range("M34").Value = -1034.4615
$ws.Range("N34").Value = -2530.3572

# Sheet CRP, row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 790.18335
$ws.Range("I58").Value = 369.69766
$ws.Range("J58").Value = 1853.7646
$ws.Range("K58").Value = 369.69766
$ws.Range("L58").Value = 1853.7646
$ws.Range("M58").Value = -166.69766
$ws.Range("N58").Value = -2259.7646

# Sheet CRP, row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 12504590
$ws.Range("I99").Value = 15629500
$ws.Range("J99").Value = 4950
$ws.Range("K99").Value = 15629500
$ws.Range("L99").Value = 4950
$ws.Range("M99").Value = -15628002
$ws.Range("N99").Value = -7946

# Sheet CRP, row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 12504590
$ws.Range("I126").Value = 15629500
$ws.Range("J126").Value = 4950
$ws.Range("K126").Value = 46888500
$ws.Range("L126").Value = 14850
$ws.Range("M126").Value = -46886030
$ws.Range("N126").Value = -19790

# Sheet CRP, row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 790.18335
$ws.Range("I136").Value = 369.69766
$ws.Range("J136").Value = 1853.7646
$ws.Range("K136").Value = 1109.09298
$ws.Range("L136").Value = 5561.293799999999
$ws.Range("M136").Value = 1440.90702
$ws.Range("N136").Value = -10661.2938

# Sheet CUL, row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 251313.16
$ws.Range("I107").Value = 276
$ws.Range("J107").Value = 401935.44
$ws.Range("K107").Value = 828
$ws.Range("L107").Value = 1205806.32
$ws.Range("M107").Value = 1092
$ws.Range("N107").Value = -1209646.32

# Sheet CUL, row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1139.55
$ws.Range("I131").Value = 409.8
$ws.Range("J131").Value = 1177.9579
$ws.Range("K131").Value = 1229.4
$ws.Range("L131").Value = 3533.8737
$ws.Range("M131").Value = 3810.6
$ws.Range("N131").Value = -13613.8737

# Sheet GSM, row 57
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 22500
$ws.Range("J57").Value = 20000
$ws.Range("L57").Value = 20000
$ws.Range("N57").Value = -21640

# Sheet GSM, row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2115.5
$ws.Range("I102").Value = 1778
$ws.Range("J102").Value = 2453
$ws.Range("K102").Value = 1778
$ws.Range("L102").Value = 2453
$ws.Range("M102").Value = -156
$ws.Range("N102").Value = -5697

# Sheet GSM, row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1855169.6
$ws.Range("I122").Value = 2781129.5
$ws.Range("K122").Value = 8343388.5
$ws.Range("M122").Value = -8340938.5

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2824.4048
$ws.Range("I132").Value = 2219.6667
$ws.Range("J132").Value = 5041.778
$ws.Range("K132").Value = 6659.000100000001
$ws.Range("L132").Value = 15125.334
$ws.Range("M132").Value = -4129.000100000001
$ws.Range("N132").Value = -20185.334

# Sheet LTW, row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 110.611115
$ws.Range("J55").Value = 145
$ws.Range("L55").Value = 145
$ws.Range("N55").Value = -491

# Sheet LTW, row 81
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H81").Value = 50000
$ws.Range("J81").Value = 50000
$ws.Range("L81").Value = 50000
$ws.Range("N81").Value = -51996

# Sheet LTW, row 84
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H84").Value = 50000
$ws.Range("J84").Value = 50000
$ws.Range("L84").Value = 150000
$ws.Range("N84").Value = -159984

# Sheet WVR, row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 58810.223
$ws.Range("I122").Value = 93380
$ws.Range("J122").Value = 4486.2856
$ws.Range("K122").Value = 280140
$ws.Range("L122").Value = 13458.8568
$ws.Range("M122").Value = -277690
$ws.Range("N122").Value = -18358.8568

# Sheet WVR, row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 167716.83
$ws.Range("I126").Value = 333800.34
$ws.Range("J126").Value = 1633.3334
$ws.Range("K126").Value = 1001401.02
$ws.Range("L126").Value = 4900.0002
$ws.Range("M126").Value = -998931.02
$ws.Range("N126").Value = -9840.0002

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 22730744
$ws.Range("I132").Value = 38464164
$ws.Range("K132").Value = 115392492
$ws.Range("M132").Value = -115389962
